$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 22
$ws1.Range("F5").Value = 2317
$ws1.Range("F6").Value = 211
$ws1.Range("F7").Value = 370

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 22
$ws4.Range("F5").Value = 2317
$ws4.Range("F6").Value = 211
$ws4.Range("F9").Value = 370
